# Update "Förändrad" (Changed) date column C for rows 2-5 from
# 2023-09-14 (serial 45183) to 2023-09-15 (serial 45184).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

foreach ($row in 2..5) {
    $ws.Cells.Item($row, 3).Value = 45184
}
